# Apply the bookmark edits described by the commit:
#  - remove the stray "_GoBack" bookmark
#  - wrap the "Структурная организация..." .. "...переработки всей
#    архитектуры." paragraphs in a new bookmark "_Hlk201141137" (id 0)
#  - wrap the "Диаграмма типа..." .. "...рисунке 2." paragraph in a new
#    bookmark "_Hlk201141164" (id 1)

$d = $word.ActiveDocument

# --- 1. Drop the old "_GoBack" bookmark (both start/end markers) -----------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- 2. New bookmark around the "component description" paragraphs --------
$r1 = $d.Content
$r1.Find.Execute("Структурная организация системы симуляции частиц") | Out-Null
$start1 = $r1.Start

$r2 = $d.Content
$r2.Find.Execute("переработки всей архитектуры.") | Out-Null
$end1 = $r2.Paragraphs(1).Range.End

$bm1 = $d.Range($start1, $end1)
$d.Bookmarks.Add("_Hlk201141137", $bm1) | Out-Null

# --- 3. New bookmark around the "entity-relationship diagram" paragraph ---
$r3 = $d.Content
$r3.Find.Execute("Диаграмма типа «сущность-связь» описывает взаимодействие") | Out-Null
$start2 = $r3.Start

$r4 = $d.Content
$r4.Find.Execute("сущностей внутри системы симуляции частиц, она приведена на рисунке 2.") | Out-Null
$end2 = $r4.Paragraphs(1).Range.End

$bm2 = $d.Range($start2, $end2)
$d.Bookmarks.Add("_Hlk201141164", $bm2) | Out-Null

Write-Output "bookmarks updated"
